# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value2 = "Datos actualizados a 7 de Agosto de 2020 a las 08:51"

# 2. Rows 145/146: Georgia and Republica de Chipre swap order, with refreshed data.
#    Row 145 becomes Georgia with new stats; Row 146 becomes Republica de Chipre
#    carrying the stats that used to belong to row 145 (unchanged totals).
$ws.Range("A145").Value2 = "Georgia"
$ws.Range("B145").Value2 = 1213
$ws.Range("C145").Value2 = 7
$ws.Range("D145").Value2 = 994
$ws.Range("E145").Value2 = 202
$ws.Range("F145").Value2 = 0
$ws.Range("G145").Value2 = 0
$ws.Range("H145").Value2 = 17

$ws.Range("A146").Value2 = "Republica de Chipre"
$ws.Range("B146").Value2 = 1208
$ws.Range("C146").Value2 = 0
$ws.Range("D146").Value2 = 856
$ws.Range("E146").Value2 = 333
$ws.Range("F146").Value2 = 0
$ws.Range("G146").Value2 = 0
$ws.Range("H146").Value2 = 19

# 3. Rows 202/203: Santa Lucia and Timor Oriental swap order (stats unchanged).
$ws.Range("A202").Value2 = "Santa Lucia"
$ws.Range("A203").Value2 = "Timor Oriental"

# 4. Refresh case counts for several countries (row data updates only).
# Row 4: Estados Unidos
$ws.Range("B4").Value2 = 5032278
$ws.Range("C4").Value2 = 99
$ws.Range("E4").Value2 = 2292806

# Row 6: India
$ws.Range("B6").Value2 = 2030001
$ws.Range("C6").Value2 = 4592
$ws.Range("D6").Value2 = 1378642
$ws.Range("E6").Value2 = 609686
$ws.Range("G6").Value2 = 35
$ws.Range("H6").Value2 = 41673

# Row 37: Ucrania
$ws.Range("B37").Value2 = 78261
$ws.Range("C37").Value2 = 1453
$ws.Range("D37").Value2 = 43055
$ws.Range("E37").Value2 = 33354
$ws.Range("G37").Value2 = 33
$ws.Range("H37").Value2 = 1852

# Row 54: Armenia
$ws.Range("B54").Value2 = 39985
$ws.Range("C54").Value2 = 166
$ws.Range("D54").Value2 = 32008
$ws.Range("E54").Value2 = 7200
$ws.Range("G54").Value2 = 5
$ws.Range("H54").Value2 = 777

# Row 62: Uzbekistan
$ws.Range("B62").Value2 = 28809
$ws.Range("C62").Value2 = 494
$ws.Range("E62").Value2 = 9043
$ws.Range("G62").Value2 = 4
$ws.Range("H62").Value2 = 179

# Row 73: El Salvador
$ws.Range("D73").Value2 = 9271
$ws.Range("E73").Value2 = 9335
$ws.Range("G73").Value2 = 7
$ws.Range("H73").Value2 = 520

# Row 141: Letonia
$ws.Range("B141").Value2 = 1281
$ws.Range("C141").Value2 = 6
$ws.Range("E141").Value2 = 179

# Row 169: Birmania
$ws.Range("B169").Value2 = 359
$ws.Range("C169").Value2 = 2
$ws.Range("E169").Value2 = 45
